$d = $word.ActiveDocument

# Replace all occurrences of the old error message prefix with the new one.
$d.Content.Find.Execute(
    "Invalid if statement: Unexpected tag EOF missing [ENDFOR]",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Invalid block: Unexpected tag EOF missing [ENDFOR]",
    2
)
